# Rename the document title casing ("Change Log" -> "Change log") and
# swap the LaTeX distro name ("LiveTeX" -> "TexLive") per the
# "name changed to Open Science Pipeline" commit.

$d = $word.ActiveDocument

# Heading1: "Change Log" -> "Change log"
$d.Content.Find.Execute(
    "Change Log", $true, $false, $false, $false, $false,
    $true, 1, $false, "Change log", 2
) | Out-Null

# Verbatim-styled tool name inside the "Added" bullet list:
# "LiveTeX" -> "TexLive"
$d.Content.Find.Execute(
    "LiveTeX", $true, $false, $false, $false, $false,
    $true, 1, $false, "TexLive", 2
) | Out-Null
